$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (A1:F1) takes on the values that were previously in row 2 (A2:F2)
$ws.Range("A1:F1").Value2 = $ws.Range("A2:F2").Value2

# Rows 2 through 11 (A2:F11) are cleared to 0
$ws.Range("A2:F11").Value2 = 0
